# Update the "time_taken" (column F) timestamps on the existing "data" sheet.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

$newTimes = @(
  "2021-10-05 14:20:43.968183",
  "2021-10-05 14:20:43.968191",
  "2021-10-05 14:20:43.968195",
  "2021-10-05 14:20:43.968197",
  "2021-10-05 14:20:43.968200",
  "2021-10-05 14:20:43.968203",
  "2021-10-05 14:20:43.968206",
  "2021-10-05 14:20:43.968208",
  "2021-10-05 14:20:43.968211",
  "2021-10-05 14:20:43.968214",
  "2021-10-05 14:20:43.968216",
  "2021-10-05 14:20:43.968219",
  "2021-10-05 14:20:43.968221",
  "2021-10-05 14:20:43.968224",
  "2021-10-05 14:20:43.968226"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet right after the "data" sheet.
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used by the "data" sheet.
$metaSheet.PageSetup.LeftMargin = $dataSheet.PageSetup.LeftMargin
$metaSheet.PageSetup.RightMargin = $dataSheet.PageSetup.RightMargin
$metaSheet.PageSetup.TopMargin = $dataSheet.PageSetup.TopMargin
$metaSheet.PageSetup.BottomMargin = $dataSheet.PageSetup.BottomMargin
$metaSheet.PageSetup.HeaderMargin = $dataSheet.PageSetup.HeaderMargin
$metaSheet.PageSetup.FooterMargin = $dataSheet.PageSetup.FooterMargin

# Match the outline summary placement used by the "data" sheet.
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Header row.
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row.
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hereditary haemorrhagic telangiectasia"
$metaSheet.Range("C2").Value = 123

# D2 must stay a text value "2.9" (not a number) while keeping the default style.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "2.9"
$dataSheet.Range("B2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)

$metaSheet.Range("E2").Value = "2021-09-27T15:14:59.476744Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:43.964449"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/123/?format=json"

# Apply the same header style (bold, centered, bordered) used on the "data" sheet.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$dataSheet.Select()
